$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row based on the sheet's used range
$lastRow = $ws.UsedRange.Rows.Count

# The workbook has columns:
#   A = code
#   B = status
#   C = codeforiati:group-name (before) / codeforiati:group-code (after)
#   D = codeforiati:group-code (before) / codeforiati:group-name (after)
# The edit swaps the contents (and header) of columns C and D for every row.
for ($r = 1; $r -le $lastRow; $r++) {
    $cVal = $ws.Cells.Item($r, 3).Value2
    $dVal = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 3).Value2 = $dVal
    $ws.Cells.Item($r, 4).Value2 = $cVal
}
